# Add more detail about dependencies in the components/dependency table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Generate sample table"): clarify the data required - it now
# depends on both global run stats and demultiplexing stats.
$ws.Range("C11").Value = "global run stats, demultiplexing stats"

# Row 12 ("Generate project email"): it also now requires demultiplexing stats.
$ws.Range("C12").Value = "demultiplexing stats"

# Row 18 ("Get run information"): the data provided is more specific now -
# it's the number of read cycles and read passes (not just "Number of reads").
$ws.Range("C18").Value = "Number of read cycles and read passes"

# New "Data required" column (G) for the supporting-code rows (15-21),
# spelling out what each piece of supporting code needs as input.
$ws.Range("G15").Value = "-"
$ws.Range("G16").Value = "-"
$ws.Range("G17").Value = "-"
$ws.Range("G18").Value = "-"
$ws.Range("G19").Value = "sample sheet, run information (num read passes)"
$ws.Range("G20").Value = "fastqc data, file names"
$ws.Range("G21").Value = "fastqc data, file names"

# Restore the selection to where the author last left it.
$null = $ws.Range("C12").Select()
